$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Pattern" column (D) values for each state row, grouped into
# the contiguous runs in which they occur (this also controls the order new
# shared strings are first created: circle, stripe, none).
$ws.Range("D2:D14").Value = "circle"
$ws.Range("D15:D28").Value = "stripe"
$ws.Range("D29:D31").Value = "circle"
$ws.Range("D32:D33").Value = "stripe"
$ws.Range("D34:D35").Value = "circle"
$ws.Range("D36:D48").Value = "none"
$ws.Range("D49").Value = "circle"

# Header for the new column, added last so "test" is appended to the shared
# string table after circle/stripe/none.
$ws.Range("D1").Value = "test"

# Select the whole new column, matching the author's last interaction.
$ws.Columns("D").Select()
